$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CH2PH")
Write-Output $ws.Name
Write-Output $ws.Range("A1").Value
